# Update "Förändrad" (column C) date for all existing data rows (2..442)
# from 2023-09-20 (45189) to 2023-09-21 (45190).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C442").Value2 = 45190

# Row 442 gains an explicit row height (ht="15" customHeight="1")
$ws.Rows.Item(442).RowHeight = 15

# Give the two new rows the same explicit row height treatment as row 443
# (row 444 keeps the default, i.e. no explicit RowHeight is set for it).
$ws.Rows.Item(443).RowHeight = 15

# --- New row 443 ---
$row = 443
$ws.Cells.Item($row, 1).Value2 = "A 44438-2023"
$ws.Cells.Item($row, 2).Value2 = 45189
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 3).Value2 = 45190
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 4).Value2 = "UPPSALA LÄN"
$ws.Cells.Item($row, 5).Value2 = "TIERP"
$ws.Cells.Item($row, 6).Value2 = "Bergvik skog öst AB"
$ws.Cells.Item($row, 7).Value2 = 6.9
$ws.Cells.Item($row, 8).Value2 = 0
$ws.Cells.Item($row, 9).Value2 = 0
$ws.Cells.Item($row, 10).Value2 = 0
$ws.Cells.Item($row, 11).Value2 = 0
$ws.Cells.Item($row, 12).Value2 = 0
$ws.Cells.Item($row, 13).Value2 = 0
$ws.Cells.Item($row, 14).Value2 = 0
$ws.Cells.Item($row, 15).Value2 = 0
$ws.Cells.Item($row, 16).Value2 = 0
$ws.Cells.Item($row, 17).Value2 = 0
$ws.Cells.Item($row, 18).Value2 = ""
$ws.Cells.Item($row, 18).WrapText = $true

# --- New row 444 ---
$row = 444
$ws.Cells.Item($row, 1).Value2 = "A 44430-2023"
$ws.Cells.Item($row, 2).Value2 = 45189
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 3).Value2 = 45190
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($row, 4).Value2 = "UPPSALA LÄN"
$ws.Cells.Item($row, 5).Value2 = "TIERP"
$ws.Cells.Item($row, 6).Value2 = "Bergvik skog öst AB"
$ws.Cells.Item($row, 7).Value2 = 1.3
$ws.Cells.Item($row, 8).Value2 = 0
$ws.Cells.Item($row, 9).Value2 = 0
$ws.Cells.Item($row, 10).Value2 = 0
$ws.Cells.Item($row, 11).Value2 = 0
$ws.Cells.Item($row, 12).Value2 = 0
$ws.Cells.Item($row, 13).Value2 = 0
$ws.Cells.Item($row, 14).Value2 = 0
$ws.Cells.Item($row, 15).Value2 = 0
$ws.Cells.Item($row, 16).Value2 = 0
$ws.Cells.Item($row, 17).Value2 = 0
$ws.Cells.Item($row, 18).Value2 = ""
$ws.Cells.Item($row, 18).WrapText = $true
